$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet used to have a blank, bold-formatted spacer row across A1:D1
# followed by the real "Name"/"Code" header row in row 2. Delete that
# spacer row so the header row becomes row 1.
$ws.Rows(1).Delete()

# The "Code" column is being dropped from the sample format - only the
# "Name" header remains, so clear out the former B1:D1 range (the rest of
# the old row extent, now occupied by the shifted-up row).
$ws.Range("B1:D1").Clear()

# Restore the recorded selection/active cell.
$ws.Range("D8").Select() | Out-Null
